$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- Sheet "About" ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: $newVersion"

$aboutWs.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Cook Colliery Coal Mine, Australia, M1531, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $dataWs.Cells.Item($dataWs.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataWs.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
